$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 05:48"

# Row 30 - Belgica
$ws.Range("B30").Value = 230480
$ws.Range("C30").Value = 8227
$ws.Range("D30").Value = 21214
$ws.Range("E30").Value = 198823
$ws.Range("G30").Value = 30
$ws.Range("H30").Value = 10443

# Row 45 - Kazajistan
$ws.Range("B45").Value = 109623
$ws.Range("C45").Value = 115
$ws.Range("D45").Value = 105145
$ws.Range("E45").Value = 2710

# Row 53 - Honduras
$ws.Range("B53").Value = 89381
$ws.Range("C53").Value = 956
$ws.Range("D53").Value = 35398
$ws.Range("E53").Value = 51407
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 2576

# Row 153 - Belice
$ws.Range("B153").Value = 2833
$ws.Range("C153").Value = 20
$ws.Range("D153").Value = 1692
$ws.Range("E153").Value = 1096
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 45

# Row 186 - Butan
$ws.Range("B186").Value = 330
$ws.Range("C186").Value = 3
$ws.Range("E186").Value = 29

# Row 187 - Mongolia
$ws.Range("B187").Value = 326
$ws.Range("C187").Value = 2
$ws.Range("E187").Value = 14
